$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.817.98"
$ws.Range("E2").Value = "  +4.15%  "
$ws.Range("D3").Value = "1.912.22"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.698"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.89%  "
$ws.Range("E9").Value = "  +5.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0759"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.816"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.45%  "
$ws.Range("D15").Value = "2.188.60"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.13%  "
$ws.Range("D17").Value = "1.904.78"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "36.815.49"
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D20").Value = "0.0₃0855"
$ws.Range("E20").Value = "  +3.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "250.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.26%  "
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.11%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0612"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0876"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +19.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +60.32%  "
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.872"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "104.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0228"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("E44").Value = "  +18.82%  "
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("D46").Value = "1.353.97"
$ws.Range("E46").Value = "  +3.86%  "
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.96%  "
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").Value = "2.090.43"
$ws.Range("E51").Value = "  +1.49%  "
